$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.735.37"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.085.63"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.51"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.98"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0779"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.17"
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.394.79"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.07"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.775"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.087.82"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.762.11"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.91"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.99"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.69"
$ws.Range("E26").Value = "  +8.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.45"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.52"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.66"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.32"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0235"
$ws.Range("E40").Value = "  +9.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.38"
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0962"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.66"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.451.11"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.07"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.20"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.278.78"
$ws.Range("E51").Value = "  +0.85%  "
